$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Shadow - true" row (row 8) is removed because cell text shadow
# formatting is not supported - only textboxes/shapes support shadows.
$ws.Rows("8").Delete()
